$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from A4 onto A5, then set its value to the new date
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 43983

# Update the active selection to B5
$ws.Range("B5").Select()
